# Update the Correspond Handoff/Handback DateTime values for the
# "614be0bb-13ad-4615-839c-1ebfb3a211fe" entries (row 4) in the
# zh-cn and de-de worksheets, as part of regenerating the handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-15 02:54:21"
$wsZhCn.Range("G4").Value = "2016-02-15 02:55:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-15 02:54:35"
$wsDeDe.Range("G4").Value = "2016-02-15 02:55:57"
